$d = $word.ActiveDocument

$d.Content.Find.Execute("今天天气不错啊", $true, $false, $false, $false, $false,
                         $true, 1, $false, "今天天气很好", 2)

$d.Content.Find.Execute("心情也很好", $true, $false, $false, $false, $false,
                         $true, 1, $false, "心情也不错", 2)
